$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns updated per data row (E,G,H,I,J,K,M,N,O,P,Q,R,S,T); F,L unchanged.
$cols = @("E","G","H","I","J","K","M","N","O","P","Q","R","S","T")

$newValues = @{
  2 = @{ "E"=3; "G"=0.6818126666666666; "H"=2.045438; "I"=0.5993607991797695; "J"=0.5993607991797694; "K"=2; "M"=2.019046; "N"=4.038092; "O"=0.003975353327590414; "P"=0.002707057536954368; "Q"=1.376611137382667; "R"=8.259666824296; "S"=0.002382670947446546; "T"=0.001622504168774588 }
  3 = @{ "E"=3; "G"=0.6818126666666666; "H"=2.045438; "I"=0.5993607991797695; "J"=0.5993607991797694; "K"=3; "M"=139.6948166666666; "N"=419.0844499999999; "O"=0.2750488370661026; "P"=0.2809459811695414; "Q"=95.24569547101109; "R"=857.2112592390998; "S"=0.1648534907974054; "T"=0.1683880078001208 }
  4 = @{ "E"=3; "G"=0.6818126666666666; "H"=2.045438; "I"=0.5993607991797695; "J"=0.5993607991797694; "K"=3; "M"=186.3548536666667; "N"=559.064561; "O"=0.3669190239530987; "P"=0.3747858972750337; "Q"=127.0590997247464; "R"=1143.531897522718; "S"=0.2199168794307902; "T"=0.2246319749120712 }
  5 = @{ "E"=3; "G"=0.6818126666666666; "H"=2.045438; "I"=0.5993607991797695; "J"=0.5993607991797694; "K"=3; "M"=143.6051993333333; "N"=430.815598; "O"=0.2827480934208787; "P"=0.2888103122968479; "Q"=97.91184390465821; "R"=881.206595141924; "S"=0.169468123239294; "T"=0.1731015795895976 }
  6 = @{ "E"=3; "G"=0.6818126666666666; "H"=2.045438; "I"=0.5993607991797695; "J"=0.5993607991797694; "K"=3; "M"=6.253715333333335; "N"=18.761146; "O"=0.01231310631861279; "P"=0.0125771036621259; "Q"=4.263862327994223; "R"=38.37476095194801; "S"=0.007379993243509229; "T"=0.007538222902298585 }
  7 = @{ "E"=3; "G"=0.6818126666666666; "H"=2.045438; "I"=0.5993607991797695; "J"=0.5993607991797694; "K"=2; "M"=29.9633245; "N"=59.926649; "O"=0.05899558591371687; "P"=0.04017364805949665; "Q"=20.42937417954366; "R"=122.576245077262; "S"=0.03535964152132409; "T"=0.0240785098069067 }
  8 = @{ "E"=3; "G"=0.4557536666666667; "H"=1.367261; "I"=0.4006392008202306; "J"=0.4006392008202306; "K"=2; "M"=2.019046; "N"=4.038092; "O"=0.003975353327590414; "P"=0.002707057536954368; "Q"=0.9201876176686666; "R"=5.521125706012; "S"=0.001592682380143868; "T"=0.00108455336817978 }
  9 = @{ "E"=3; "G"=0.4557536666666667; "H"=1.367261; "I"=0.4006392008202306; "J"=0.4006392008202306; "K"=3; "M"=139.6948166666666; "N"=419.0844499999999; "O"=0.2750488370661026; "P"=0.2809459811695414; "Q"=63.6664249101611; "R"=572.9978241914499; "S"=0.1101953462686971; "T"=0.1125579733694206 }
  10 = @{ "E"=3; "G"=0.4557536666666667; "H"=1.367261; "I"=0.4006392008202306; "J"=0.4006392008202306; "K"=3; "M"=186.3548536666667; "N"=559.064561; "O"=0.3669190239530987; "P"=0.3747858972750337; "Q"=84.93190785971345; "R"=764.3871707374211; "S"=0.1470021445223085; "T"=0.1501539223629625 }
  11 = @{ "E"=3; "G"=0.4557536666666667; "H"=1.367261; "I"=0.4006392008202306; "J"=0.4006392008202306; "K"=3; "M"=143.6051993333333; "N"=430.815598; "O"=0.2827480934208787; "P"=0.2888103122968479; "Q"=65.44859614856422; "R"=589.0373653370781; "S"=0.1132799701815847; "T"=0.1157087327072504 }
  12 = @{ "E"=3; "G"=0.4557536666666667; "H"=1.367261; "I"=0.4006392008202306; "J"=0.4006392008202306; "K"=3; "M"=6.253715333333335; "N"=18.761146; "O"=0.01231310631861279; "P"=0.0125771036621259; "Q"=2.850153693456223; "R"=25.65138324110601; "S"=0.004933113075103558; "T"=0.005038880759827317 }
  13 = @{ "E"=3; "G"=0.4557536666666667; "H"=1.367261; "I"=0.4006392008202306; "J"=0.4006392008202306; "K"=2; "M"=29.9633245; "N"=59.926649; "O"=0.05899558591371687; "P"=0.04017364805949665; "Q"=13.65589500639817; "R"=81.935370038389; "S"=0.02363594439239278; "T"=0.01609513825258994 }
}

foreach ($row in $newValues.Keys) {
  foreach ($col in $cols) {
    $ws.Range("$col$row").Value = $newValues[$row][$col]
  }
}
